# ---------------------------------------------------------------
# PyNeapple_test_results.xlsx - "Testing with new parameters"
# Updates raw measurement data on the "absolut" sheet; the
# "relativ" sheet recomputes its D:G/I:J formulas automatically
# on recalculation (run_com recalcs after the script completes).
# ---------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("absolut")

# Row 4
$ws.Range("D4").Value = [double]"2.9583665618994931E-3"
$ws.Range("E4").Value = [double]"0.4314881750747997"
$ws.Range("F4").Value = [double]"1.9165124293185299E-3"
$ws.Range("G4").Value = [double]"0.22086502669545791"

# Row 5
$ws.Range("D5").Value = [double]"8.2415177630682911E-3"
$ws.Range("E5").Value = [double]"0.54980550062417044"
$ws.Range("F5").Value = [double]"6.626775671996335E-3"
$ws.Range("G5").Value = [double]"0.89140640312641017"

# Row 6
$ws.Range("D6").Value = [double]"0.130353828973573"
$ws.Range("E6").Value = [double]"4.7879841012607859E-2"
$ws.Range("F6").Value = [double]"0.14466833238158541"
$ws.Range("G6").Value = [double]"5.0959051860008361E-2"

# Row 7
$ws.Range("D7").Value = [double]"0.1910026074596459"
$ws.Range("E7").Value = [double]"3.475941023606554E-2"

# Row 8
$ws.Range("D8").Value = [double]"2.9583665618994931E-3"
$ws.Range("E8").Value = [double]"0.4314881750747997"
$ws.Range("F8").Value = [double]"1.9165124293185299E-3"
$ws.Range("G8").Value = [double]"0.22086502669545791"

# Row 9
$ws.Range("D9").Value = [double]"8.2415177630682911E-3"
$ws.Range("E9").Value = [double]"0.54980550062417044"
$ws.Range("F9").Value = [double]"6.626775671996335E-3"
$ws.Range("G9").Value = [double]"0.89140640312641017"

# Row 10
$ws.Range("D10").Value = [double]"0.130353828973573"
$ws.Range("E10").Value = [double]"4.7879841012607859E-2"
$ws.Range("F10").Value = [double]"0.14466833238158541"
$ws.Range("G10").Value = [double]"5.0959051860008361E-2"

# Row 11
$ws.Range("D11").Value = [double]"0.1910026074596459"
$ws.Range("E11").Value = [double]"3.475941023606554E-2"

# Row 13
$ws.Range("D13").Value = [double]"8.4733626044853011E-4"
$ws.Range("E13").Value = [double]"0.14690428482304824"
$ws.Range("F13").Value = [double]"8.4733626044853011E-4"
$ws.Range("G13").Value = [double]"0.14526614972371041"

# Row 14
$ws.Range("D14").Value = [double]"6.2422531657229056E-3"
$ws.Range("E14").Value = [double]"0.7567534537955708"
$ws.Range("F14").Value = [double]"6.1357258277597148E-3"
$ws.Range("G14").Value = [double]"0.76411398617563442"

# Row 15
$ws.Range("D15").Value = [double]"0.26566231827347758"
$ws.Range("E15").Value = [double]"9.3273023901056004E-2"
$ws.Range("F15").Value = [double]"0.26566231827347758"
$ws.Range("G15").Value = [double]"9.2515079156380478E-2"

# Row 16
$ws.Range("D16").Value = [double]"0.26566231827347758"
$ws.Range("E16").Value = [double]"9.465329387860727E-2"

# Row 22
$ws.Range("D22").Value = [double]"1.61099985825943E-3"
$ws.Range("E22").Value = [double]"0.70243247259634778"
$ws.Range("F22").Value = [double]"1.1187205415065359E-3"
$ws.Range("G22").Value = [double]"0.62605681373561139"

# Row 23
$ws.Range("D23").Value = [double]"9.635665677864837E-3"
$ws.Range("E23").Value = [double]"0.22545232964853451"
$ws.Range("F23").Value = [double]"6.0291786234855493E-3"
$ws.Range("G23").Value = [double]"0.39744101704231144"

# Row 24
$ws.Range("D24").Value = [double]"0.1497809921024926"
$ws.Range("E24").Value = [double]"9.8309419491559719E-2"
$ws.Range("F24").Value = [double]"0.15240468466226989"
$ws.Range("G24").Value = [double]"0.12008145558233596"

# Row 25
$ws.Range("D25").Value = [double]"0.20474143893734029"
$ws.Range("E25").Value = [double]"8.6955321020121823E-2"

# Row 26
$ws.Range("D26").Value = [double]"1.61099985825943E-3"
$ws.Range("E26").Value = [double]"0.70243247259634778"
$ws.Range("F26").Value = [double]"1.1187205415065359E-3"
$ws.Range("G26").Value = [double]"0.62605681373561139"

# Row 27
$ws.Range("D27").Value = [double]"9.635665677864837E-3"
$ws.Range("E27").Value = [double]"0.22545232964853451"
$ws.Range("F27").Value = [double]"6.0291786234855493E-3"
$ws.Range("G27").Value = [double]"0.39744101704231144"

# Row 28
$ws.Range("D28").Value = [double]"0.1497809921024926"
$ws.Range("E28").Value = [double]"9.8309419491559719E-2"
$ws.Range("F28").Value = [double]"0.15240468466226989"
$ws.Range("G28").Value = [double]"0.12008145558233596"

# Row 29
$ws.Range("D29").Value = [double]"0.20474143893734029"
$ws.Range("E29").Value = [double]"8.6955321020121823E-2"

# Row 30
$ws.Range("D30").Value = [double]"0.21197711581854251"
$ws.Range("E30").Value = [double]"6.8839707399211775E-2"

# Row 32
$ws.Range("D32").Value = [double]"8.7728159746556708E-4"
$ws.Range("E32").Value = [double]"0.58947376395780982"
$ws.Range("F32").Value = [double]"8.7728159746556698E-4"
$ws.Range("G32").Value = [double]"0.58961735864715192"

# Row 33
$ws.Range("D33").Value = [double]"5.2471712739034426E-3"
$ws.Range("E33").Value = [double]"0.28926239543313736"
$ws.Range("F33").Value = [double]"5.2471712739034426E-3"
$ws.Range("G33").Value = [double]"0.28926239543313736"

# Row 34
$ws.Range("D34").Value = [double]"0.2156902891245093"
$ws.Range("E34").Value = [double]"0.11611746062601215"
$ws.Range("F34").Value = [double]"0.2156902891245093"
$ws.Range("G34").Value = [double]"0.11686327110461661"

# Row 35
$ws.Range("D35").Value = [double]"0.26108886143117699"
$ws.Range("E35").Value = [double]"0.10304111241787529"

# Row 12 / Row 31 trailing values cleared (group now has one fewer data point)
$ws.Range("D12").ClearContents()
$ws.Range("E12").ClearContents()
$ws.Range("D31").ClearContents()
$ws.Range("E31").ClearContents()

# "NNLS_AUC" / "reg = 0" row now marked with an asterisk, plus a new
# "*" annotation column (O:Q) duplicating the updated reg=0 series
$ws.Range("B22").Value = "reg = 0*"
$ws.Range("O22").Value = "*"
$ws.Range("P22").Value = [double]"1.61099985825943E-3"
$ws.Range("Q22").Value = [double]"0.70243247259634778"
$ws.Range("P23").Value = [double]"9.635665677864837E-3"
$ws.Range("Q23").Value = [double]"0.22545232964853451"
$ws.Range("P24").Value = [double]"0.1497809921024926"
$ws.Range("Q24").Value = [double]"9.8309419491559719E-2"
$ws.Range("P25").Value = [double]"0.20474143893734029"
$ws.Range("Q25").Value = [double]"8.6955321020121823E-2"
$ws.Range("P26").Value = [double]"0.21197711581854251"
$ws.Range("Q26").Value = [double]"6.8839707399211775E-2"

# Restore the original selections / active sheet so the saved view
# matches (relativ stays on J12, absolut becomes the active tab at B22)
$ws1 = $wb.Worksheets.Item("relativ")
$ws1.Range("J12").Select()
$ws.Activate()
$ws.Range("B22").Select()
